# issue #5: stock data output to json file
# Add a "property_category" column to the 股票 (stock) sheet, between the
# existing "total" and "date" columns, and populate it with "stock" for
# every existing data row. This shifts the old date/owner(legislator_name)/
# legislator_id columns one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Determine the extent of the existing data (header row + data rows).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Insert a new column at H, shifting the subsequent columns (date,
# legislator_name, legislator_id) one column to the right.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Every existing record on this sheet describes a stock holding.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
